$wb = $excel.ActiveWorkbook

# Rename the second worksheet from "具有相當價值之財產" to "保險"
$ws = $wb.Worksheets.Item("具有相當價值之財產")
$ws.Name = "保險"

# Delete the last row (row 6) so the table shrinks from 6 to 5 data rows
$ws.Rows.Item(6).Delete()

# Rewrite the remaining rows with the new insurance data
# Row 1 (was the header row) now carries the first policy's data, no headers
$ws.Cells.Item(1, 1).Value = ""
$ws.Cells.Item(1, 2).Value = "國泰人壽"
$ws.Cells.Item(1, 3).Value = "保本111終身"
$ws.Cells.Item(1, 4).Value = "陳淑慧"
$ws.Cells.Item(1, 5).Value = ""

# Row 2
$ws.Cells.Item(2, 1).Value = 101
$ws.Cells.Item(2, 2).Value = "國泰人壽"
$ws.Cells.Item(2, 3).Value = "保本111終身"
$ws.Cells.Item(2, 4).Value = "陳淑慧"
$ws.Cells.Item(2, 5).Value = ""

# Row 3
$ws.Cells.Item(3, 1).Value = 102
$ws.Cells.Item(3, 2).Value = "國泰人壽"
$ws.Cells.Item(3, 3).Value = "美意年年终生壽險"
$ws.Cells.Item(3, 4).Value = "陳淑慧"
$ws.Cells.Item(3, 5).Value = ""

# Row 4
$ws.Cells.Item(4, 1).Value = 103
$ws.Cells.Item(4, 2).Value = "紐約人壽"
$ws.Cells.Item(4, 3).Value = "聚寶盆變額萬能壽險"
$ws.Cells.Item(4, 4).Value = "陳淑慧"
$ws.Cells.Item(4, 5).Value = "躉繳"

# Row 5
$ws.Cells.Item(5, 1).Value = 104
$ws.Cells.Item(5, 2).Value = "保誠人尋"
$ws.Cells.Item(5, 3).Value = "築夢人生"
$ws.Cells.Item(5, 4).Value = "陳淑慧"
$ws.Cells.Item(5, 5).Value = "躉繳"
